$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Split the "BGBSC" sheet: the old 2019 column (B) is dropped, the old
#    2020 column (C) becomes the sole data column of a brand-new "SYGBSC"
#    sheet placed right before "BGBSC", and the remaining years (old D:AG,
#    2021-2050) shift left to become the new B:AE range on "BGBSC".
# ---------------------------------------------------------------------------

$bgbsc = $wb.Worksheets.Item("BGBSC")

# Remember formatting / label info we will need once the new sheet exists.
$tabColor = $bgbsc.Tab.Color
$rowLabel = $bgbsc.Range("A2").Value()

# Capture the "2020" column (currently column C) before any shifting, since
# it is about to move onto the new sheet.
$year2020 = $bgbsc.Range("C1").Value()
$formula2020 = $bgbsc.Range("C2").Formula
$numFmt2020 = $bgbsc.Range("C2").NumberFormat

# Drop the old 2019 column (B) completely - shifts 2020..2050 left to B:AF.
$bgbsc.Range("B1:B2").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# Create the new sheet right before BGBSC and name it.
$sygbsc = $wb.Worksheets.Add($bgbsc)
$sygbsc.Name = "SYGBSC"

# Re-fetch BGBSC - the handle above becomes stale once a sheet is inserted.
$bgbsc = $wb.Worksheets.Item("BGBSC")

# Populate SYGBSC with the single 2020 data column.
$sygbsc.Range("A2").Value = $rowLabel
$sygbsc.Range("B1").Value = $year2020
$sygbsc.Range("B2").Formula = $formula2020
$sygbsc.Range("B2").NumberFormat = $numFmt2020
$sygbsc.Tab.Color = $tabColor

# Finally drop the (now first) 2020 column from BGBSC so the remaining
# years 2021-2050 shift left into B:AE.
$bgbsc.Range("B1:B2").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# ---------------------------------------------------------------------------
# 2) Turn on iterative calculation (matches calcPr iterate="1"
#    iterateDelta="1.0000000000000001E-5" in the target workbook).
# ---------------------------------------------------------------------------
$wb.Iteration = $true
$wb.MaxChange = 0.00001
